# Data.xlsx - "Add files via upload" commit
# Changes on the "SA" sheet:
#  1. Rename three codes (same Name (Roman)/Name (Original) text, just the
#     short Code value shifts):
#       B13: SAPT -> SAPR   (Shiva Prarthana)
#       B14: SADP -> SAPT   (Shri Adishaktyah Prarthana)
#       B15: SAPS -> SAPW   (Shri Annapurneshwari Stotram)
#  2. Insert a new row after "Shri Ganesha Atharvashirsha" (old row 23),
#     pushing everything below down by one row, with:
#       Code            = SGAU
#       Name (Roman)    = Shri Ganesha Atharvashirsha (Complete)
#       Name (Original) = श्री गणेश अथर्वशीर्ष (संपूर्ण)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SA")
$ws.Activate()

# 1. Code renames (rows 13-15)
$ws.Range("B13").Value = "SAPR"
$ws.Range("B14").Value = "SAPT"
$ws.Range("B15").Value = "SAPW"

# 2. Insert a new row at 24 (shifts old row 24.. down to 25..)
$ws.Rows.Item(24).Insert()

# Re-establish the running serial-number formula through the newly
# inserted row (Excel's relative-fill semantics keep the A2:A44 count
# sequential: 1,2,3 ... 43).
$ws.Range("A24:A25").Formula = "=A23+1"

$ws.Range("B24").Value = "SGAU"
$ws.Range("C24").Value = "Shri Ganesha Atharvashirsha (Complete)"
$ws.Range("D24").Value = "श्री गणेश अथर्वशीर्ष (संपूर्ण)"

# Match the author's final selection / scroll position on the sheet.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select()
